$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.261.66'
$ws.Range('E2').Value = '  +0.92%  '
$ws.Range('D3').Value = '1.653.10'
$ws.Range('E3').Value = '  +0.49%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').NumberFormat = 'General'
$ws.Range('E4').Value = '  -0.27%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '218.90'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('E5').Value = '  -0.39%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.511'
$ws.Range('D6').NumberFormat = 'General'
$ws.Range('E6').Value = '  +2.25%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('D7').NumberFormat = 'General'
$ws.Range('E7').Value = '  -0.30%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.257'
$ws.Range('D8').NumberFormat = 'General'
$ws.Range('E8').Value = '  +1.46%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0630'
$ws.Range('D9').NumberFormat = 'General'
$ws.Range('E9').Value = '  +0.85%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.28'
$ws.Range('D10').NumberFormat = 'General'
$ws.Range('E10').Value = '  +4.30%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0848'
$ws.Range('D11').NumberFormat = 'General'
$ws.Range('E11').Value = '  +0.08%  '
$ws.Range('D12').Value = '1.884.20'
$ws.Range('E12').Value = '  +0.50%  '
$ws.Range('D13').Value = '1.646.59'
$ws.Range('E13').Value = '  +0.14%  '
$ws.Range('E14').Value = '  -0.52%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.541'
$ws.Range('D15').NumberFormat = 'General'
$ws.Range('E15').Value = '  +2.08%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '67.90'
$ws.Range('D16').NumberFormat = 'General'
$ws.Range('E16').Value = '  +2.84%  '
$ws.Range('D17').Value = '27.256.54'
$ws.Range('E17').Value = '  +1.00%  '
$ws.Range('E18').Value = '  +0.85%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '221.26'
$ws.Range('D19').NumberFormat = 'General'
$ws.Range('E19').Value = '  +1.24%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.00'
$ws.Range('D20').NumberFormat = 'General'
$ws.Range('E20').Value = '  -0.38%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.79'
$ws.Range('D21').NumberFormat = 'General'
$ws.Range('E21').Value = '  +2.50%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.45'
$ws.Range('D22').NumberFormat = 'General'
$ws.Range('E22').Value = '  +1.07%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.51'
$ws.Range('D23').NumberFormat = 'General'
$ws.Range('E23').Value = '  +2.69%  '
$ws.Range('E24').Value = '  +0.59%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '148.07'
$ws.Range('D25').NumberFormat = 'General'
$ws.Range('E25').Value = '  -0.15%  '
$ws.Range('E26').Value = '  -0.31%  '
$ws.Range('E27').Value = '  +1.14%  '
$ws.Range('E28').Value = '  +1.10%  '
$ws.Range('E29').Value = '  +0.38%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0509'
$ws.Range('D30').NumberFormat = 'General'
$ws.Range('E30').Value = '  -0.51%  '
$ws.Range('E31').Value = '  -0.38%  '
$ws.Range('E32').Value = '  -0.13%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.03'
$ws.Range('D33').NumberFormat = 'General'
$ws.Range('E33').Value = '  +1.03%  '
$ws.Range('E34').Value = '  +2.16%  '
$ws.Range('D35').Value = '1.274.90'
$ws.Range('E35').Value = '  +0.32%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.45'
$ws.Range('D36').NumberFormat = 'General'
$ws.Range('E36').Value = '  +0.73%  '
$ws.Range('E37').Value = '  +3.22%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.543'
$ws.Range('D38').NumberFormat = 'General'
$ws.Range('E38').Value = '  +1.92%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.847'
$ws.Range('D39').NumberFormat = 'General'
$ws.Range('E39').Value = '  +2.64%  '
$ws.Range('E40').Value = '  -0.28%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.815'
$ws.Range('D41').NumberFormat = 'General'
$ws.Range('E41').Value = '  +1.00%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.41'
$ws.Range('D42').NumberFormat = 'General'
$ws.Range('E42').Value = '  +0.84%  '
$ws.Range('B43').Value = 'RocketPoolETH'
$ws.Range('C43').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D43').Value = '1.794.27'
$ws.Range('E43').Value = '  +0.53%  '
$ws.Range('B44').Value = 'MXToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.17'
$ws.Range('D44').NumberFormat = 'General'
$ws.Range('E44').Value = '  +5.96%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '63.22'
$ws.Range('D45').NumberFormat = 'General'
$ws.Range('E45').Value = '  +2.19%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '92.74'
$ws.Range('D46').NumberFormat = 'General'
$ws.Range('E46').Value = '  +0.23%  '
$ws.Range('E47').Value = '  +0.00%  '
$ws.Range('E48').Value = '  +16.02%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.72'
$ws.Range('D50').NumberFormat = 'General'
$ws.Range('E50').Value = '  +1.61%  '
$ws.Range('E51').Value = '  +0.96%  '
